$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 271; this shifts the existing rows 271-405
# down to 272-406, preserving all of their original values/styles.
$ws.Rows(271).Insert()

# Populate the newly inserted row 271 with the new record.
$ws.Range("A271").Value = 5
$ws.Range("B271").Value = "Macroferia Regional de Talca"
$ws.Range("C271").Value = "Maule"
$ws.Range("D271").Value2 = 45016
$ws.Range("D271").NumberFormat = $ws.Range("D272").NumberFormat
$ws.Range("E271").Value = 7
$ws.Range("F271").Value = 100112009
$ws.Range("G271").Value = "Acelga"
$ws.Range("H271").Value = "Sin especificar"
$ws.Range("I271").Value = "Primera"
$ws.Range("J271").Value = 500
$ws.Range("K271").Value = 2500
$ws.Range("L271").Value = 2500
$ws.Range("M271").Value = 2500
$ws.Range("N271").Value = "$/docena de atados (4 kilos)"
$ws.Range("O271").Value = "Región del Maule"
$ws.Range("P271").Value = 625
$ws.Range("Q271").Value = 4
$ws.Range("R271").Value = "Hortaliza"
